$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for phase columns
$ws.Range("K3").Value = "phase 3"
$ws.Range("L3").Value = "phase 4"
$ws.Range("M3").Value = "phase 5"

# Row 4: 5% loss data
$ws.Range("J4").Value = "5% loss"
$ws.Range("L4").Value = 15743
$ws.Range("M4").Value = 17336

# Row 5: 10%error data
$ws.Range("J5").Value = "10%error"
$ws.Range("L5").Value = 28697
$ws.Range("M5").Value = 29538

# Update selection to match target state
$ws.Range("K4").Select()
